$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, pushing the existing 5001 "Bullet" row down to row 9.
$ws.Rows("8:8").Insert()

# Drop the stale row-outline grouping level left over on the sheet (no row actually
# uses row grouping any more).
$ws.Rows("8:9").Ungroup()

# Copy the formatting (styles) of row 4 (the first data row, same style pattern we need)
# onto the new row 8, so B8:J8 pick up styles 3,3,3,3,3,7,3,3,3.
$ws.Range("B4:J4").Copy()
$ws.Range("B8:J8").PasteSpecial(-4122)

# Row 4 has no A cell either; make sure the new row 8 doesn't have one (the Insert()
# carried one down from row 7's formatting).
$ws.Range("A8").Clear()

# Fill in the data for the new unit (2001 / Food).
$ws.Range("B8").Value = 2001
$ws.Range("C8").Value = "Food"
$ws.Range("D8").Value = "测试说明"
$ws.Range("E8").Value = "UnitConfig.Name.2001"
$ws.Range("F8").Value = "米克尔"
$ws.Range("G8").Value = "带有强力攻击技能"
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 178
$ws.Range("J8").Value = 68

$ws.Range("E8").Select() | Out-Null
